$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 100, pushing the existing rows 100.. down to 102..
$ws.Rows.Item(100).Resize(2).Insert()

# Row 100 - new "Primera" quality record
$ws.Cells.Item(100, 1).Value = 4
$ws.Cells.Item(100, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(100, 3).Value = 'Los Lagos'
$ws.Cells.Item(100, 4).Value = 44603
$ws.Cells.Item(100, 5).Value = 10
$ws.Cells.Item(100, 6).Value = 'Fruta'
$ws.Cells.Item(100, 7).Value = 100108
$ws.Cells.Item(100, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(100, 9).Value = 100108002
$ws.Cells.Item(100, 10).Value = 'Mango'
$ws.Cells.Item(100, 11).Value = 'Sin especificar'
$ws.Cells.Item(100, 12).Value = 'Primera'
$ws.Cells.Item(100, 13).Value = 160
$ws.Cells.Item(100, 14).Value = 7500
$ws.Cells.Item(100, 15).Value = 8000
$ws.Cells.Item(100, 16).Value = 7750
$ws.Cells.Item(100, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(100, 18).Value = 'Perú'
$ws.Cells.Item(100, 19).Value = 1938
$ws.Cells.Item(100, 20).Value = 4

# Row 101 - new "Segunda" quality record (same fecha/market as row 100)
$ws.Cells.Item(101, 1).Value = 4
$ws.Cells.Item(101, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(101, 3).Value = 'Los Lagos'
$ws.Cells.Item(101, 4).Value = 44603
$ws.Cells.Item(101, 5).Value = 10
$ws.Cells.Item(101, 6).Value = 'Fruta'
$ws.Cells.Item(101, 7).Value = 100108
$ws.Cells.Item(101, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(101, 9).Value = 100108002
$ws.Cells.Item(101, 10).Value = 'Mango'
$ws.Cells.Item(101, 11).Value = 'Sin especificar'
$ws.Cells.Item(101, 12).Value = 'Segunda'
$ws.Cells.Item(101, 13).Value = 60
$ws.Cells.Item(101, 14).Value = 5000
$ws.Cells.Item(101, 15).Value = 5000
$ws.Cells.Item(101, 16).Value = 5000
$ws.Cells.Item(101, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(101, 18).Value = 'Perú'
$ws.Cells.Item(101, 19).Value = 1250
$ws.Cells.Item(101, 20).Value = 4
